# Auto-generated script to append scrim result rows per the commit diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: Kaboom Canyon ---
$ws = $wb.Worksheets.Item("Kaboom Canyon")
# Row 59
$ws.Range("A58:N58").Copy($ws.Range("A59:N59"))
$ws.Range("A59").Value = "CHUCK"
$ws.Range("B59").Value = "FINX"
$ws.Range("C59").Value = "BELLE"
$ws.Range("D59").Value = "BONNIE"
$ws.Range("E59").Value = "ANGELO"
$ws.Range("F59").Value = "MEG"
$ws.Range("G59").Value = "Equipo 1"
$ws.Range("H59").Value = "nyamura"
$ws.Range("I59").Value = "RC|Shu"
$ws.Range("J59").Value = "RC|Battoman"
$ws.Range("K59").Value = "ZETA|Levi"
$ws.Range("L59").Value = "ZETA|Sizuku"
$ws.Range("M59").Value = "ZETA|Sitetampo"
$ws.Range("N59").Value = "20250710T141520.000Z"
$ws.Range("G59").Interior.Color = 16770508

# Row 60
$ws.Range("A58:N58").Copy($ws.Range("A60:N60"))
$ws.Range("A60").Value = "CHUCK"
$ws.Range("B60").Value = "FINX"
$ws.Range("C60").Value = "BELLE"
$ws.Range("D60").Value = "BONNIE"
$ws.Range("E60").Value = "ANGELO"
$ws.Range("F60").Value = "MEG"
$ws.Range("G60").Value = "Equipo 1"
$ws.Range("H60").Value = "nyamura"
$ws.Range("I60").Value = "RC|Shu"
$ws.Range("J60").Value = "RC|Battoman"
$ws.Range("K60").Value = "ZETA|Levi"
$ws.Range("L60").Value = "ZETA|Sizuku"
$ws.Range("M60").Value = "ZETA|Sitetampo"
$ws.Range("N60").Value = "20250710T141325.000Z"
$ws.Range("G60").Interior.Color = 16770508


# --- Sheet: Pit Stop ---
$ws = $wb.Worksheets.Item("Pit Stop")
# Row 89
$ws.Range("A88:N88").Copy($ws.Range("A89:N89"))
$ws.Range("A89").Value = "DRACO"
$ws.Range("B89").Value = "CHARLIE"
$ws.Range("C89").Value = "KAZE"
$ws.Range("D89").Value = "HANK"
$ws.Range("E89").Value = "MICO"
$ws.Range("F89").Value = "LUMI"
$ws.Range("G89").Value = "Equipo 2"
$ws.Range("H89").Value = "HMB|BosS"
$ws.Range("I89").Value = "HMB|Lukii"
$ws.Range("J89").Value = "HMB|Symantec"
$ws.Range("K89").Value = "TH|Code: LeNain"
$ws.Range("L89").Value = "TH|IKAUSSA"
$ws.Range("M89").Value = "TH|Zhar"
$ws.Range("N89").Value = "20250710T141355.000Z"
$ws.Range("G89").Interior.Color = 13421812

# Row 90
$ws.Range("A88:N88").Copy($ws.Range("A90:N90"))
$ws.Range("A90").Value = "DRACO"
$ws.Range("B90").Value = "CHARLIE"
$ws.Range("C90").Value = "KAZE"
$ws.Range("D90").Value = "HANK"
$ws.Range("E90").Value = "MICO"
$ws.Range("F90").Value = "LUMI"
$ws.Range("G90").Value = "Equipo 1"
$ws.Range("H90").Value = "HMB|BosS"
$ws.Range("I90").Value = "HMB|Lukii"
$ws.Range("J90").Value = "HMB|Symantec"
$ws.Range("K90").Value = "TH|Code: LeNain"
$ws.Range("L90").Value = "TH|IKAUSSA"
$ws.Range("M90").Value = "TH|Zhar"
$ws.Range("N90").Value = "20250710T141116.000Z"
$ws.Range("G90").Interior.Color = 16770508

# Row 91
$ws.Range("A88:N88").Copy($ws.Range("A91:N91"))
$ws.Range("A91").Value = "DRACO"
$ws.Range("B91").Value = "CHARLIE"
$ws.Range("C91").Value = "KAZE"
$ws.Range("D91").Value = "HANK"
$ws.Range("E91").Value = "MICO"
$ws.Range("F91").Value = "LUMI"
$ws.Range("G91").Value = "Equipo 2"
$ws.Range("H91").Value = "HMB|BosS"
$ws.Range("I91").Value = "HMB|Lukii"
$ws.Range("J91").Value = "HMB|Symantec"
$ws.Range("K91").Value = "TH|Code: LeNain"
$ws.Range("L91").Value = "TH|IKAUSSA"
$ws.Range("M91").Value = "TH|Zhar"
$ws.Range("N91").Value = "20250710T140908.000Z"
$ws.Range("G91").Interior.Color = 13421812


# --- Sheet: Goldarm Gulch ---
$ws = $wb.Worksheets.Item("Goldarm Gulch")
# Row 70
$ws.Range("A69:N69").Copy($ws.Range("A70:N70"))
$ws.Range("A70").Value = "BELLE"
$ws.Range("B70").Value = "OLLIE"
$ws.Range("C70").Value = "BROCK"
$ws.Range("D70").Value = "HANK"
$ws.Range("E70").Value = "GUS"
$ws.Range("F70").Value = "PIPER"
$ws.Range("G70").Value = "Equipo 1"
$ws.Range("H70").Value = "nyamura"
$ws.Range("I70").Value = "RC|Battoman"
$ws.Range("J70").Value = "RC|Shu"
$ws.Range("K70").Value = "ZETA|Levi"
$ws.Range("L70").Value = "ZETA|Sitetampo"
$ws.Range("M70").Value = "ZETA|Sizuku"
$ws.Range("N70").Value = "20250710T140645.000Z"
$ws.Range("G70").Interior.Color = 16770508

# Row 71
$ws.Range("A69:N69").Copy($ws.Range("A71:N71"))
$ws.Range("A71").Value = "BELLE"
$ws.Range("B71").Value = "OLLIE"
$ws.Range("C71").Value = "BROCK"
$ws.Range("D71").Value = "HANK"
$ws.Range("E71").Value = "GUS"
$ws.Range("F71").Value = "PIPER"
$ws.Range("G71").Value = "Equipo 2"
$ws.Range("H71").Value = "nyamura"
$ws.Range("I71").Value = "RC|Battoman"
$ws.Range("J71").Value = "RC|Shu"
$ws.Range("K71").Value = "ZETA|Levi"
$ws.Range("L71").Value = "ZETA|Sitetampo"
$ws.Range("M71").Value = "ZETA|Sizuku"
$ws.Range("N71").Value = "20250710T140426.000Z"
$ws.Range("G71").Interior.Color = 13421812

# Row 72
$ws.Range("A69:N69").Copy($ws.Range("A72:N72"))
$ws.Range("A72").Value = "BONNIE"
$ws.Range("B72").Value = "GUS"
$ws.Range("C72").Value = "CHARLIE"
$ws.Range("D72").Value = "KAZE"
$ws.Range("E72").Value = "BELLE"
$ws.Range("F72").Value = "BROCK"
$ws.Range("G72").Value = "Equipo 1"
$ws.Range("H72").Value = "CR|Moya"
$ws.Range("I72").Value = "CR|Milkreo"
$ws.Range("J72").Value = "Tensai 천재"
$ws.Range("K72").Value = "NAVI|Ryohei"
$ws.Range("L72").Value = "NAVI|Kuru"
$ws.Range("M72").Value = "NAVI|Achapi"
$ws.Range("N72").Value = "20250710T141348.000Z"
$ws.Range("G72").Interior.Color = 16770508

# Row 73
$ws.Range("A69:N69").Copy($ws.Range("A73:N73"))
$ws.Range("A73").Value = "LUMI"
$ws.Range("B73").Value = "BROCK"
$ws.Range("C73").Value = "KAZE"
$ws.Range("D73").Value = "HANK"
$ws.Range("E73").Value = "WILLOW"
$ws.Range("F73").Value = "JAE-YONG"
$ws.Range("G73").Value = "Equipo 1"
$ws.Range("H73").Value = "CR|Moya"
$ws.Range("I73").Value = "CR|Milkreo"
$ws.Range("J73").Value = "Tensai 천재"
$ws.Range("K73").Value = "NAVI|Ryohei"
$ws.Range("L73").Value = "NAVI|Achapi"
$ws.Range("M73").Value = "NAVI|Kuru"
$ws.Range("N73").Value = "20250710T140815.000Z"
$ws.Range("G73").Interior.Color = 16770508

# Row 74
$ws.Range("A69:N69").Copy($ws.Range("A74:N74"))
$ws.Range("A74").Value = "LUMI"
$ws.Range("B74").Value = "BROCK"
$ws.Range("C74").Value = "KAZE"
$ws.Range("D74").Value = "HANK"
$ws.Range("E74").Value = "WILLOW"
$ws.Range("F74").Value = "JAE-YONG"
$ws.Range("G74").Value = "Equipo 2"
$ws.Range("H74").Value = "CR|Moya"
$ws.Range("I74").Value = "CR|Milkreo"
$ws.Range("J74").Value = "Tensai 천재"
$ws.Range("K74").Value = "NAVI|Ryohei"
$ws.Range("L74").Value = "NAVI|Achapi"
$ws.Range("M74").Value = "NAVI|Kuru"
$ws.Range("N74").Value = "20250710T140458.000Z"
$ws.Range("G74").Interior.Color = 13421812


# --- Sheet: Open Business ---
$ws = $wb.Worksheets.Item("Open Business")
# Row 50
$ws.Range("A49:N49").Copy($ws.Range("A50:N50"))
$ws.Range("A50").Value = "HANK"
$ws.Range("B50").Value = "BEA"
$ws.Range("C50").Value = "ALLI"
$ws.Range("D50").Value = "LUMI"
$ws.Range("E50").Value = "DRACO"
$ws.Range("F50").Value = "SHADE"
$ws.Range("G50").Value = "Equipo 1"
$ws.Range("H50").Value = "GEN|Moding"
$ws.Range("I50").Value = "GEN|cookie"
$ws.Range("J50").Value = "GEN|BONOX2"
$ws.Range("K50").Value = "RVL|Terry"
$ws.Range("L50").Value = "Isee in Bush"
$ws.Range("M50").Value = "RVL|Mameshi"
$ws.Range("N50").Value = "20250710T140441.000Z"
$ws.Range("G50").Interior.Color = 16770508

# Row 51
$ws.Range("A49:N49").Copy($ws.Range("A51:N51"))
$ws.Range("A51").Value = "HANK"
$ws.Range("B51").Value = "BEA"
$ws.Range("C51").Value = "ALLI"
$ws.Range("D51").Value = "LUMI"
$ws.Range("E51").Value = "DRACO"
$ws.Range("F51").Value = "SHADE"
$ws.Range("G51").Value = "Equipo 1"
$ws.Range("H51").Value = "GEN|Moding"
$ws.Range("I51").Value = "GEN|cookie"
$ws.Range("J51").Value = "GEN|BONOX2"
$ws.Range("K51").Value = "RVL|Terry"
$ws.Range("L51").Value = "Isee in Bush"
$ws.Range("M51").Value = "RVL|Mameshi"
$ws.Range("N51").Value = "20250710T140231.000Z"
$ws.Range("G51").Interior.Color = 16770508


# --- Sheet: Ring of Fire ---
$ws = $wb.Worksheets.Item("Ring of Fire")
# Row 68
$ws.Range("A67:N67").Copy($ws.Range("A68:N68"))
$ws.Range("A68").Value = "FINX"
$ws.Range("B68").Value = "BROCK"
$ws.Range("C68").Value = "AMBER"
$ws.Range("D68").Value = "JUJU"
$ws.Range("E68").Value = "LUMI"
$ws.Range("F68").Value = "CROW"
$ws.Range("G68").Value = "Equipo 1"
$ws.Range("H68").Value = "GEN|cookie"
$ws.Range("I68").Value = "GEN|BONOX2"
$ws.Range("J68").Value = "GEN|Moding"
$ws.Range("K68").Value = "Isee in Bush"
$ws.Range("L68").Value = "RVL|Terry"
$ws.Range("M68").Value = "RVL|Mameshi"
$ws.Range("N68").Value = "20250710T141357.000Z"
$ws.Range("G68").Interior.Color = 16770508

# Row 69
$ws.Range("A67:N67").Copy($ws.Range("A69:N69"))
$ws.Range("A69").Value = "FINX"
$ws.Range("B69").Value = "BROCK"
$ws.Range("C69").Value = "AMBER"
$ws.Range("D69").Value = "JUJU"
$ws.Range("E69").Value = "LUMI"
$ws.Range("F69").Value = "CROW"
$ws.Range("G69").Value = "Equipo 2"
$ws.Range("H69").Value = "GEN|cookie"
$ws.Range("I69").Value = "GEN|BONOX2"
$ws.Range("J69").Value = "GEN|Moding"
$ws.Range("K69").Value = "Isee in Bush"
$ws.Range("L69").Value = "RVL|Terry"
$ws.Range("M69").Value = "RVL|Mameshi"
$ws.Range("N69").Value = "20250710T141113.000Z"
$ws.Range("G69").Interior.Color = 13421812


# --- Sheet: Hideout ---
$ws = $wb.Worksheets.Item("Hideout")
# Row 52
$ws.Range("A51:N51").Copy($ws.Range("A52:N52"))
$ws.Range("A52").Value = "BELLE"
$ws.Range("B52").Value = "DRACO"
$ws.Range("C52").Value = "JAE-YONG"
$ws.Range("D52").Value = "DOUG"
$ws.Range("E52").Value = "LUMI"
$ws.Range("F52").Value = "GRAY"
$ws.Range("G52").Value = "Equipo 1"
$ws.Range("H52").Value = "KDS|Decaii"
$ws.Range("I52").Value = "KDS|Remica"
$ws.Range("J52").Value = "KDS|Ćiro"
$ws.Range("K52").Value = "NOVO|26is"
$ws.Range("L52").Value = "NOVO|Biso"
$ws.Range("M52").Value = "NOVO|Marco"
$ws.Range("N52").Value = "20250710T141236.000Z"
$ws.Range("G52").Interior.Color = 16770508

# Row 53
$ws.Range("A51:N51").Copy($ws.Range("A53:N53"))
$ws.Range("A53").Value = "BELLE"
$ws.Range("B53").Value = "DRACO"
$ws.Range("C53").Value = "JAE-YONG"
$ws.Range("D53").Value = "DOUG"
$ws.Range("E53").Value = "LUMI"
$ws.Range("F53").Value = "GRAY"
$ws.Range("G53").Value = "Equipo 1"
$ws.Range("H53").Value = "KDS|Decaii"
$ws.Range("I53").Value = "KDS|Remica"
$ws.Range("J53").Value = "KDS|Ćiro"
$ws.Range("K53").Value = "NOVO|26is"
$ws.Range("L53").Value = "NOVO|Biso"
$ws.Range("M53").Value = "NOVO|Marco"
$ws.Range("N53").Value = "20250710T141016.000Z"
$ws.Range("G53").Interior.Color = 16770508

